$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the cells we are about to overwrite as Text so the
# numeric-looking strings (prices / percentages) are preserved
# verbatim instead of being re-interpreted as numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.37%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.06"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.07%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.121"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.29%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07852"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-2.08%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.266"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-10.32%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.838"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.21%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.809"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.60%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9279"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1770"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.17%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07821"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "7.02%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08897"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.26%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03086"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.13%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1004"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.08%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001515"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.04%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005877"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.27%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.459"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.89%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.249"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.09%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.23%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1329"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.14%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.252"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "25.42%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1794"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.72%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04602"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.29%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001252"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.68%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004500"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.00%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.10%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-1.36%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01797"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.33%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04836"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "8.29%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007241"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.35%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1375"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.41%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002121"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.16%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009947"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "3.75%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006270"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.68%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.01%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003597"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-58.90%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7770"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-5.31%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.01%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.01%"
